# Rename the five sheets to describe the number of latent classes modelled
# on each, fix up the chart series formulas (Excel does NOT do this for us
# automatically when a sheet is renamed via COM), drop the leftover
# "_xlchart.v1.*" helper defined names that Excel had stashed for the
# chart-wizard caches, restore the remaining selections, and move the
# active tab from "Five Classes" to "Six Classes".

$wb = $excel.ActiveWorkbook

$renames = @{
    "Sheet1" = "Two Classes"
    "Sheet2" = "Three Classes"
    "Sheet3" = "Four Classes"
    "Sheet4" = "Five Classes"
    "Sheet5" = "Six Classes"
}

# 1. Fix up every chart series formula BEFORE renaming so the old sheet
#    name is still what is stored in each c:f reference; simple text
#    substitution mirrors what the author's click-through rename would
#    have produced once Excel regenerated the chart cache.
foreach ($old in $renames.Keys) {
    $new = $renames[$old]
    $ws = $wb.Worksheets.Item($old)
    foreach ($co in $ws.ChartObjects()) {
        $chart = $co.Chart
        for ($i = 1; $i -le $chart.SeriesCollection().Count; $i++) {
            $series = $chart.SeriesCollection().Item($i)
            $series.Formula = $series.Formula.Replace("$old!", "'$new'!")
        }
    }
}

# 2. Rename the worksheets themselves - this also fixes up the
#    defined names (_xlnm._FilterDatabase) that reference them.
foreach ($old in $renames.Keys) {
    $wb.Worksheets.Item($old).Name = $renames[$old]
}

# 3. Drop the unused "_xlchart.v1.*" defined names left behind by the
#    chart wizard - they are no longer needed.
for ($i = 0; $i -le 6; $i++) {
    $wb.Names.Item("_xlchart.v1.$i").Delete()
}

# 4. Restore per-sheet selections that moved around during the edit.
$wb.Worksheets.Item("Two Classes").Range("B22").Select()
$wb.Worksheets.Item("Three Classes").Range("B10").Select()

# 5. Make "Six Classes" the active tab (was "Five Classes" before).
$wb.Worksheets.Item("Six Classes").Activate()
$wb.Worksheets.Item("Six Classes").Range("C19").Select()
